$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 93: take data (B:AC) from original row 96 (id/A column stays put)
$ws.Range("B93").Value = 6236612
$ws.Range("F93").Value = 'Zamora'
$ws.Range("G93").Value = 'Carabobo'
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 2
$ws.Range("J93").Value = 'A'
$ws.Range("K93").Value = 3.2
$ws.Range("L93").Value = 3.1
$ws.Range("M93").Value = 2.15
$ws.Range("N93").Value = 4.5
$ws.Range("O93").Value = 3.3
$ws.Range("P93").Value = 1.75
$ws.Range("Q93").Value = 0.5
$ws.Range("R93").Value = 2
$ws.Range("S93").Value = 1.8
$ws.Range("T93").Value = 2.25
$ws.Range("U93").Value = 1.925
$ws.Range("V93").Value = 1.875
$ws.Range("W93").Value = -1
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = 0.75
$ws.Range("Z93").Value = -1
$ws.Range("AA93").Value = 0.8
$ws.Range("AB93").Value = -0.5
$ws.Range("AC93").Value = 0.4375

# Row 94: take data (B:AC) from original row 97 (id/A column stays put)
$ws.Range("B94").Value = 6236252
$ws.Range("F94").Value = 'Deportivo Tachira'
$ws.Range("G94").Value = 'CD Hermanos Colmenares'
$ws.Range("H94").Value = 1
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 'H'
$ws.Range("K94").Value = 1.363
$ws.Range("L94").Value = 4.2
$ws.Range("M94").Value = 7.5
$ws.Range("N94").Value = 1.333
$ws.Range("O94").Value = 4.5
$ws.Range("P94").Value = 8
$ws.Range("Q94").Value = -1.5
$ws.Range("R94").Value = 2
$ws.Range("S94").Value = 1.8
$ws.Range("T94").Value = 2.5
$ws.Range("U94").Value = 1.925
$ws.Range("V94").Value = 1.875
$ws.Range("W94").Value = 0.333
$ws.Range("X94").Value = -1
$ws.Range("Y94").Value = -1
$ws.Range("Z94").Value = -1
$ws.Range("AA94").Value = 0.8
$ws.Range("AB94").Value = -1
$ws.Range("AC94").Value = 0.875

# Row 95: take data (B:AC) from original row 98 (id/A column stays put)
$ws.Range("B95").Value = 6236611
$ws.Range("F95").Value = 'Mineros'
$ws.Range("G95").Value = 'Monagas'
$ws.Range("H95").Value = 2
$ws.Range("I95").Value = 1
$ws.Range("J95").Value = 'H'
$ws.Range("K95").Value = 3.2
$ws.Range("L95").Value = 3.4
$ws.Range("M95").Value = 2
$ws.Range("N95").Value = 4.2
$ws.Range("O95").Value = 3.8
$ws.Range("P95").Value = 1.65
$ws.Range("Q95").Value = 0.75
$ws.Range("R95").Value = 1.95
$ws.Range("S95").Value = 1.85
$ws.Range("T95").Value = 2.5
$ws.Range("U95").Value = 1.825
$ws.Range("V95").Value = 1.975
$ws.Range("W95").Value = 3.2
$ws.Range("X95").Value = -1
$ws.Range("Y95").Value = -1
$ws.Range("Z95").Value = 0.95
$ws.Range("AA95").Value = -1
$ws.Range("AB95").Value = 0.825
$ws.Range("AC95").Value = -1

# Row 96: take data (B:AC) from original row 94 (id/A column stays put)
$ws.Range("B96").Value = 6236255
$ws.Range("F96").Value = 'Deportivo Rayo Zuliano'
$ws.Range("G96").Value = 'Caracas'
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 'D'
$ws.Range("K96").Value = 3.75
$ws.Range("L96").Value = 3.1
$ws.Range("M96").Value = 1.95
$ws.Range("N96").Value = 2.9
$ws.Range("O96").Value = 2.875
$ws.Range("P96").Value = 2.45
$ws.Range("Q96").Value = 0.25
$ws.Range("R96").Value = 1.775
$ws.Range("S96").Value = 2.025
$ws.Range("T96").Value = 2.25
$ws.Range("U96").Value = 1.85
$ws.Range("V96").Value = 1.95
$ws.Range("W96").Value = -1
$ws.Range("X96").Value = 1.875
$ws.Range("Y96").Value = -1
$ws.Range("Z96").Value = 0.3875
$ws.Range("AA96").Value = -0.5
$ws.Range("AB96").Value = -1
$ws.Range("AC96").Value = 0.95

# Row 97: take data (B:AC) from original row 99 (id/A column stays put)
$ws.Range("B97").Value = 6236253
$ws.Range("F97").Value = 'Deportivo La Guaira'
$ws.Range("G97").Value = 'UCV'
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 'D'
$ws.Range("K97").Value = 1.833
$ws.Range("L97").Value = 3.25
$ws.Range("M97").Value = 4
$ws.Range("N97").Value = 2
$ws.Range("O97").Value = 3.2
$ws.Range("P97").Value = 3.5
$ws.Range("Q97").Value = -0.25
$ws.Range("R97").Value = 1.775
$ws.Range("S97").Value = 2.025
$ws.Range("T97").Value = 2.25
$ws.Range("U97").Value = 1.9
$ws.Range("V97").Value = 1.9
$ws.Range("W97").Value = -1
$ws.Range("X97").Value = 2.2
$ws.Range("Y97").Value = -1
$ws.Range("Z97").Value = -0.5
$ws.Range("AA97").Value = 0.5125
$ws.Range("AB97").Value = -1
$ws.Range("AC97").Value = 0.8999999999999999

# Row 98: take data (B:AC) from original row 93 (id/A column stays put)
$ws.Range("B98").Value = 6236254
$ws.Range("F98").Value = 'Academia Puerto Cabello'
$ws.Range("G98").Value = 'Estudiantes Merida'
$ws.Range("H98").Value = 1
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 'H'
$ws.Range("K98").Value = 1.727
$ws.Range("L98").Value = 3.4
$ws.Range("M98").Value = 4.333
$ws.Range("N98").Value = 1.666
$ws.Range("O98").Value = 3.4
$ws.Range("P98").Value = 4.75
$ws.Range("Q98").Value = -0.75
$ws.Range("R98").Value = 1.875
$ws.Range("S98").Value = 1.925
$ws.Range("T98").Value = 2.5
$ws.Range("U98").Value = 1.9
$ws.Range("V98").Value = 1.9
$ws.Range("W98").Value = 0.6659999999999999
$ws.Range("X98").Value = -1
$ws.Range("Y98").Value = -1
$ws.Range("Z98").Value = 0.4375
$ws.Range("AA98").Value = -0.5
$ws.Range("AB98").Value = -1
$ws.Range("AC98").Value = 0.8999999999999999

# Row 99: take data (B:AC) from original row 95 (id/A column stays put)
$ws.Range("B99").Value = 6236251
$ws.Range("F99").Value = 'Angostura FC'
$ws.Range("G99").Value = 'Portuguesa'
$ws.Range("H99").Value = 1
$ws.Range("I99").Value = 2
$ws.Range("J99").Value = 'A'
$ws.Range("K99").Value = 3.1
$ws.Range("L99").Value = 3.2
$ws.Range("M99").Value = 2.15
$ws.Range("N99").Value = 4
$ws.Range("O99").Value = 3.6
$ws.Range("P99").Value = 1.75
$ws.Range("Q99").Value = 0.75
$ws.Range("R99").Value = 1.8
$ws.Range("S99").Value = 2
$ws.Range("T99").Value = 2.5
$ws.Range("U99").Value = 1.95
$ws.Range("V99").Value = 1.85
$ws.Range("W99").Value = -1
$ws.Range("X99").Value = -1
$ws.Range("Y99").Value = 0.75
$ws.Range("Z99").Value = -0.5
$ws.Range("AA99").Value = 0.5
$ws.Range("AB99").Value = 0.95
$ws.Range("AC99").Value = -1

# Row 114: take data (B:AC) from original row 115 (id/A column stays put)
$ws.Range("B114").Value = 7352251
$ws.Range("F114").Value = 'Caracas'
$ws.Range("G114").Value = 'Academia Puerto Cabello'
$ws.Range("H114").Value = 1
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 'H'
$ws.Range("K114").Value = 2.1
$ws.Range("L114").Value = 3.2
$ws.Range("M114").Value = 3.3
$ws.Range("N114").Value = 2.15
$ws.Range("O114").Value = 3.1
$ws.Range("P114").Value = 3.2
$ws.Range("Q114").Value = -0.5
$ws.Range("R114").Value = 2.025
$ws.Range("S114").Value = 1.775
$ws.Range("T114").Value = 2.25
$ws.Range("U114").Value = 1.975
$ws.Range("V114").Value = 1.825
$ws.Range("W114").Value = 1.15
$ws.Range("X114").Value = -1
$ws.Range("Y114").Value = -1
$ws.Range("Z114").Value = 1.025
$ws.Range("AA114").Value = -1
$ws.Range("AB114").Value = -1
$ws.Range("AC114").Value = 0.825

# Row 115: take data (B:AC) from original row 114 (id/A column stays put)
$ws.Range("B115").Value = 7352250
$ws.Range("F115").Value = 'Portuguesa'
$ws.Range("G115").Value = 'Deportivo Tachira'
$ws.Range("H115").Value = 1
$ws.Range("I115").Value = 1
$ws.Range("J115").Value = 'D'
$ws.Range("K115").Value = 3.1
$ws.Range("L115").Value = 2.875
$ws.Range("M115").Value = 2.3
$ws.Range("N115").Value = 3
$ws.Range("O115").Value = 2.875
$ws.Range("P115").Value = 2.375
$ws.Range("Q115").Value = 0.25
$ws.Range("R115").Value = 1.725
$ws.Range("S115").Value = 2.075
$ws.Range("T115").Value = 2
$ws.Range("U115").Value = 1.825
$ws.Range("V115").Value = 1.975
$ws.Range("W115").Value = -1
$ws.Range("X115").Value = 1.875
$ws.Range("Y115").Value = -1
$ws.Range("Z115").Value = 0.3625
$ws.Range("AA115").Value = -0.5
$ws.Range("AB115").Value = 0
$ws.Range("AC115").Value = -0

# Row 116: take data (B:AC) from original row 117 (id/A column stays put)
$ws.Range("B116").Value = 7352252
$ws.Range("F116").Value = 'Deportivo Tachira'
$ws.Range("G116").Value = 'Caracas'
$ws.Range("H116").Value = 1
$ws.Range("I116").Value = 1
$ws.Range("J116").Value = 'D'
$ws.Range("K116").Value = 2.3
$ws.Range("L116").Value = 2.875
$ws.Range("M116").Value = 3.1
$ws.Range("N116").Value = 2.25
$ws.Range("O116").Value = 2.8
$ws.Range("P116").Value = 3.25
$ws.Range("Q116").Value = -0.25
$ws.Range("R116").Value = 1.975
$ws.Range("S116").Value = 1.825
$ws.Range("T116").Value = 2
$ws.Range("U116").Value = 1.925
$ws.Range("V116").Value = 1.875
$ws.Range("W116").Value = -1
$ws.Range("X116").Value = 1.8
$ws.Range("Y116").Value = -1
$ws.Range("Z116").Value = -0.5
$ws.Range("AA116").Value = 0.4125
$ws.Range("AB116").Value = 0
$ws.Range("AC116").Value = -0

# Row 117: take data (B:AC) from original row 116 (id/A column stays put)
$ws.Range("B117").Value = 7352254
$ws.Range("F117").Value = 'Academia Puerto Cabello'
$ws.Range("G117").Value = 'Portuguesa'
$ws.Range("H117").Value = 1
$ws.Range("I117").Value = 1
$ws.Range("J117").Value = 'D'
$ws.Range("K117").Value = 2.05
$ws.Range("L117").Value = 3.4
$ws.Range("M117").Value = 3
$ws.Range("N117").Value = 1.833
$ws.Range("O117").Value = 3.5
$ws.Range("P117").Value = 3.5
$ws.Range("Q117").Value = -0.25
$ws.Range("R117").Value = 1.65
$ws.Range("S117").Value = 2.2
$ws.Range("T117").Value = 2.25
$ws.Range("U117").Value = 1.825
$ws.Range("V117").Value = 1.975
$ws.Range("W117").Value = -1
$ws.Range("X117").Value = 2.5
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = -0.5
$ws.Range("AA117").Value = 0.6000000000000001
$ws.Range("AB117").Value = -0.5
$ws.Range("AC117").Value = 0.4875

# Row 162: take data (B:AC) from original row 163 (id/A column stays put)
$ws.Range("B162").Value = 7952893
$ws.Range("F162").Value = 'UCV'
$ws.Range("G162").Value = 'Deportivo La Guaira'
$ws.Range("H162").Value = 1
$ws.Range("I162").Value = 1
$ws.Range("J162").Value = 'D'
$ws.Range("K162").Value = 2.1
$ws.Range("L162").Value = 3
$ws.Range("M162").Value = 3.25
$ws.Range("N162").Value = 2.25
$ws.Range("O162").Value = 3.1
$ws.Range("P162").Value = 2.9
$ws.Range("Q162").Value = -0.25
$ws.Range("R162").Value = 2.025
$ws.Range("S162").Value = 1.775
$ws.Range("T162").Value = 2
$ws.Range("U162").Value = 1.8
$ws.Range("V162").Value = 2
$ws.Range("W162").Value = -1
$ws.Range("X162").Value = 2.1
$ws.Range("Y162").Value = -1
$ws.Range("Z162").Value = -0.5
$ws.Range("AA162").Value = 0.3875
$ws.Range("AB162").Value = 0
$ws.Range("AC162").Value = -0

# Row 163: take data (B:AC) from original row 162 (id/A column stays put)
$ws.Range("B163").Value = 7952905
$ws.Range("F163").Value = 'Angostura FC'
$ws.Range("G163").Value = 'Deportivo Tachira'
$ws.Range("H163").Value = 2
$ws.Range("I163").Value = 0
$ws.Range("J163").Value = 'H'
$ws.Range("K163").Value = 3.6
$ws.Range("L163").Value = 3.6
$ws.Range("M163").Value = 1.8
$ws.Range("N163").Value = 3.75
$ws.Range("O163").Value = 2.875
$ws.Range("P163").Value = 2.1
$ws.Range("Q163").Value = 0.25
$ws.Range("R163").Value = 1.95
$ws.Range("S163").Value = 1.85
$ws.Range("T163").Value = 2
$ws.Range("U163").Value = 2.025
$ws.Range("V163").Value = 1.775
$ws.Range("W163").Value = 2.75
$ws.Range("X163").Value = -1
$ws.Range("Y163").Value = -1
$ws.Range("Z163").Value = 0.95
$ws.Range("AA163").Value = -1
$ws.Range("AB163").Value = 0
$ws.Range("AC163").Value = -0
